$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain-text values (e.g. "64.606.85",
# "1.00") rather than numbers. Any of the new prices below that Excel
# would otherwise auto-parse as a genuine number need the cell format
# switched to Text first so the literal string is preserved.
$textPriceCells = "D5","D6","D7","D9","D10","D11","D12","D15","D16","D19","D20","D21","D22","D23","D24","D25","D26","D27","D30","D31","D34","D36","D38","D39","D40","D41","D43","D44","D45","D46","D47","D49","D50","D51"
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume(1h) / coin / link values.
$ws.Range("D2").Value = '64.606.85'
$ws.Range("E2").Value = '  +3.05%  '
$ws.Range("D3").Value = '3.466.19'
$ws.Range("E3").Value = '  +3.99%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '578.52'
$ws.Range("E5").Value = '  +4.28%  '
$ws.Range("D6").Value = '157.94'
$ws.Range("E6").Value = '  +4.36%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.467.06'
$ws.Range("E8").Value = '  +3.94%  '
$ws.Range("D9").Value = '0.559'
$ws.Range("E9").Value = '  +5.88%  '
$ws.Range("D10").Value = '7.58'
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +6.48%  '
$ws.Range("D12").Value = '0.449'
$ws.Range("E12").Value = '  +3.40%  '
$ws.Range("D13").Value = '4.067.26'
$ws.Range("E13").Value = '  +4.15%  '
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("D15").Value = '0.0000198'
$ws.Range("E15").Value = '  +9.48%  '
$ws.Range("D16").Value = '27.89'
$ws.Range("E16").Value = '  +3.74%  '
$ws.Range("D17").Value = '64.644.66'
$ws.Range("E17").Value = '  +3.08%  '
$ws.Range("D18").Value = '3.468.83'
$ws.Range("E18").Value = '  +4.14%  '
$ws.Range("D19").Value = '6.45'
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("D20").Value = '14.42'
$ws.Range("E20").Value = '  +4.88%  '
$ws.Range("D21").Value = '397.10'
$ws.Range("E21").Value = '  +2.24%  '
$ws.Range("D22").Value = '8.55'
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("D23").Value = '0.548'
$ws.Range("E23").Value = '  +1.87%  '
$ws.Range("D24").Value = '73.26'
$ws.Range("E24").Value = '  +3.35%  '
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.48%  '
$ws.Range("D26").Value = '0.0000123'
$ws.Range("E26").Value = '  +26.99%  '
$ws.Range("D27").Value = '9.63'
$ws.Range("E27").Value = '  +7.84%  '
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '6.07'
$ws.Range("E30").Value = '  +8.89%  '
$ws.Range("D31").Value = '6.72'
$ws.Range("E31").Value = '  +4.54%  '
$ws.Range("E32").Value = '  +7.21%  '
$ws.Range("E33").Value = '  +3.59%  '
$ws.Range("D34").Value = '23.86'
$ws.Range("E34").Value = '  +3.91%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").Value = '7.05'
$ws.Range("E36").Value = '  +4.63%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").Value = '160.52'
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").Value = '0.0791'
$ws.Range("E39").Value = '  +7.17%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '27.68'
$ws.Range("E40").Value = '  +1.63%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '1.88'
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").Value = '2.921.61'
$ws.Range("E42").Value = '  +2.16%  '
$ws.Range("D43").Value = '0.0325'
$ws.Range("E43").Value = '  +3.55%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.778'
$ws.Range("E44").Value = '  +3.35%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = '4.45'
$ws.Range("E45").Value = '  +2.74%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '42.30'
$ws.Range("E46").Value = '  +4.06%  '
$ws.Range("D47").Value = '23.92'
$ws.Range("E47").Value = '  +8.72%  '
$ws.Range("E48").Value = '  +5.30%  '
$ws.Range("D49").Value = '2.22'
$ws.Range("E49").Value = '  +25.38%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '6.58'
$ws.Range("E50").Value = '  +4.42%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").Value = '0.858'
$ws.Range("E51").Value = '  +6.54%  '
